$p = $ppt.ActivePresentation

# Slide 14 ("Energy level of top songs") is a leftover duplicate slide
# (a near-copy of slide 13 with a different image set) that was left in
# the deck by mistake. Remove it; PowerPoint will renumber/relink the
# remaining slides automatically on save.
$p.Slides.Item(14).Delete()
